# RR calculation.xlsx - "added lowest stop loss percentage"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Row 3: remove the hard-coded E3 value so F3 (=E2/E3) turns into a
#    #DIV/0! error, matching the target workbook.
# ---------------------------------------------------------------------
$ws.Range("E3").ClearContents()

# ---------------------------------------------------------------------
# 2. Row 6: "lowest stop loss" percentage changes from 0.15% to 0.25%,
#    and its display format switches from the built-in 0.00% to a
#    custom 0.0% format (numFmt 165). Fill (yellow) stays the same.
# ---------------------------------------------------------------------
$ws.Range("D6").Value = 0.0025
$ws.Range("D6").NumberFormat = "0.0%"

# ---------------------------------------------------------------------
# 3. New "commission" label + Win/Lose/Net mini table (M4:P9, L6).
#    Shared strings are introduced in this order: commission, Win,
#    Lose, Net -- matching the order they appear in the target file.
# ---------------------------------------------------------------------
$ws.Range("L6").Value = "commission"

$ws.Range("M4").Value = "Win"
$ws.Range("N4").Value = "Lose"
$ws.Range("P4").Value = "Net"

# Header row (M4:P4) bold + centered.
$hdr = $ws.Range("M4:P4")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108   # xlCenter

# Orange fill highlight for Win / Lose / Net (not the blank spacer O4).
$ws.Range("M4:N4").Interior.Color = 49407   # RGB(255,192,0)
$ws.Range("P4").Interior.Color = 49407

# Row 5: Win = H2 (reward), Lose = E2 (risk amount)
$ws.Range("M5").Formula = "=H2"
$ws.Range("N5").Formula = "=E2"
$ws.Range("M5:P5").Style = "Comma"

# Row 6: commission values
$ws.Range("M6").Formula = "=E6*(0.00027*2)"
$ws.Range("N6").Formula = "=M6"
$ws.Range("M6:P6").Style = "Comma"

# Row 7: totals with a top thin / bottom medium border underneath
$ws.Range("M7").Formula = "=M5-M6"
$ws.Range("N7").Formula = "=SUM(N5:N6)"
$ws.Range("P7").Formula = "=M7-N7"
$ws.Range("O7").Style = "Comma"

$totals = $ws.Range("M7:P7")
$totals.NumberFormat = "_(* #,##0.00_);_(* \(#,##0.00\);_(* ""-""??_);_(@_)"
$topB = $totals.Borders.Item(8)    # xlEdgeTop
$topB.LineStyle = 1
$topB.Weight = 2                   # xlThin
$botB = $totals.Borders.Item(9)    # xlEdgeBottom
$botB.LineStyle = 1
$botB.Weight = -4138               # xlMedium
$ws.Rows.Item(7).RowHeight = 15.75

# Row 8: Net result as a percentage of the rolling capital
$ws.Range("P8").Formula = "=P7/E2"
$ws.Range("P8").NumberFormat = "_(* #,##0.00_);_(* \(#,##0.00\);_(* ""-""??_);_(@_)"

# Row 9: stray Comma-formatted helper cell
$ws.Range("M9").NumberFormat = "_(* #,##0.00_);_(* \(#,##0.00\);_(* ""-""??_);_(@_)"

# ---------------------------------------------------------------------
# 4. Re-apply the plain Comma format to the G/I helper columns of the
#    two existing risk blocks (rows 6 & 10) -- unaffected numerically,
#    only the underlying style bookkeeping shifts because style 9 (the
#    old D6 percent+fill combo) was removed above.
# ---------------------------------------------------------------------
foreach ($addr in @("G6", "I6", "G10", "I10")) {
    $ws.Range($addr).NumberFormat = "_(* #,##0.00_);_(* \(#,##0.00\);_(* ""-""??_);_(@_)"
}

# ---------------------------------------------------------------------
# 5. New helper block in row 16/17.
# ---------------------------------------------------------------------
$ws.Range("E16").Formula = "=460*0.0004"
$ws.Range("F16").Formula = "=E16*2"
foreach ($addr in @("E16", "F16")) {
    $ws.Range($addr).NumberFormat = "_(* #,##0.00_);_(* \(#,##0.00\);_(* ""-""??_);_(@_)"
}
$ws.Range("G16").Style = "Comma"

$ws.Range("G17").NumberFormat = "0.0000%"
$ws.Range("H17").NumberFormat = "0.00000%"
$ws.Range("G17:H17").Style = "Percent"
$ws.Range("G17").NumberFormat = "0.0000%"
$ws.Range("H17").NumberFormat = "0.00000%"

# ---------------------------------------------------------------------
# 6. Column widths.
# ---------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 11.14     # -> stored width 12
$ws.Columns.Item(8).ColumnWidth = 11.14     # -> stored width 12 (was 11.14)
$ws.Columns.Item(12).ColumnWidth = 10.65    # -> stored width ~11.57

# ---------------------------------------------------------------------
# 7. Page setup + selection.
# ---------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9          # xlPaperA4
$ws.PageSetup.Orientation = 1        # xlPortrait

$ws.Range("K17").Select()
